# Adds a new row for the "lsT" (Aggregate labor supply) variable to the
# "3dim" worksheet, and makes that sheet the active/selected tab, matching
# the commit "3 dimension sim and year on same plot".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("3dim")

# Fill in the new row (row 6), entering values in the order that reproduces
# the shared-string insertion order seen in the target workbook:
#   variable_label (B), year_example (G), variable_name (A), chart_type (F), export flag (I)
$ws.Range("B6").Value = "Aggregate labor supply"
$ws.Range("G6").Value = "2025, 2030"
$ws.Range("A6").Value = "lsT"
$ws.Range("F6").Value = "sim_%bau_bar"
$ws.Range("I6").Value = 1

# Select the new row's first cell, then activate the sheet so it becomes
# both the selected range and the active tab in the workbook.
$ws.Range("A6").Select()
$ws.Activate()
